$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 202.66667
$ws.Range("I33").Value = 216.75
$ws.Range("K33").Value = 216.75
$ws.Range("M33").Value = 12.25
# Row 86
$ws.Range("H86").Value = 2787.9375
$ws.Range("I86").Value = 964.64703
$ws.Range("K86").Value = 964.64703
$ws.Range("M86").Value = 158.35297
# Row 89
$ws.Range("H89").Value = 2787.9375
$ws.Range("I89").Value = 964.64703
$ws.Range("K89").Value = 4823.23515
$ws.Range("M89").Value = 792.7648500000005
# Row 106
$ws.Range("H106").Value = 2427.2856
$ws.Range("I106").Value = 1998.5
$ws.Range("K106").Value = 1998.5
$ws.Range("M106").Value = -1367.5
# Row 107
$ws.Range("H107").Value = 717.53845
$ws.Range("I107").Value = 312.9
$ws.Range("J107").Value = 2066.3333
$ws.Range("K107").Value = 312.9
$ws.Range("L107").Value = 2066.3333
$ws.Range("M107").Value = 1607.1
$ws.Range("N107").Value = -5906.3333
# Row 132
$ws.Range("H132").Value = 31829.145
$ws.Range("I132").Value = 32250.92
$ws.Range("K132").Value = 96752.75999999999
$ws.Range("M132").Value = -94222.75999999999
# Row 137
$ws.Range("I137").Value = 1389938.4
$ws.Range("J137").Value = 1018687.75
$ws.Range("K137").Value = 4169815.2
$ws.Range("L137").Value = 3056063.25
$ws.Range("M137").Value = -4167265.2
$ws.Range("N137").Value = -3061163.25
# Row 140
$ws.Range("H140").Value = 189999.5
$ws.Range("J140").Value = 189999.5
$ws.Range("L140").Value = 189999.5
$ws.Range("N140").Value = -200359.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5052490.5
$ws.Range("I32").Value = 5052490.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5052490.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -5052203.5
$ws.Range("N32").ClearContents()
# Row 61
$ws.Range("H61").Value = 836309
$ws.Range("I61").Value = 1043841.7
$ws.Range("J61").Value = 6178.25
$ws.Range("K61").Value = 1043841.7
$ws.Range("L61").Value = 6178.25
$ws.Range("M61").Value = -1043629.7
$ws.Range("N61").Value = -6602.25
# Row 132
$ws.Range("H132").Value = 467582.1
$ws.Range("J132").Value = 10100
$ws.Range("L132").Value = 30300
$ws.Range("N132").Value = -35360
# Row 136
$ws.Range("H136").Value = 836309
$ws.Range("I136").Value = 1043841.7
$ws.Range("J136").Value = 6178.25
$ws.Range("K136").Value = 3131525.1
$ws.Range("L136").Value = 18534.75
$ws.Range("M136").Value = -3128975.1
$ws.Range("N136").Value = -23634.75

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 49000
$ws.Range("J20").Value = 49000
$ws.Range("L20").Value = 49000
$ws.Range("N20").Value = -49472
# Row 30
$ws.Range("H30").Value = 49000
$ws.Range("J30").Value = 49000
$ws.Range("L30").Value = 49000
$ws.Range("N30").Value = -49182
# Row 31
$ws.Range("H31").Value = 3904.6292
$ws.Range("I31").Value = 2192.5667
$ws.Range("K31").Value = 2192.5667
$ws.Range("M31").Value = -1897.5667
# Row 34
$ws.Range("H34").Value = 3904.6292
$ws.Range("I34").Value = 2192.5667
$ws.Range("K34").Value = 2192.5667
$ws.Range("M34").Value = -1990.5667
# Row 41
$ws.Range("H41").Value = 34533.582
$ws.Range("I41").Value = 19618.25
$ws.Range("J41").Value = 41991.25
$ws.Range("K41").Value = 19618.25
$ws.Range("L41").Value = 41991.25
$ws.Range("M41").Value = -19190.25
$ws.Range("N41").Value = -42847.25
# Row 107
$ws.Range("H107").Value = 1930.7778
$ws.Range("I107").Value = 1510.6666
$ws.Range("J107").Value = 2771
$ws.Range("K107").Value = 1510.6666
$ws.Range("L107").Value = 2771
$ws.Range("M107").Value = 409.3334
$ws.Range("N107").Value = -6611
# Row 122
$ws.Range("H122").Value = 3139.65
$ws.Range("I122").Value = 1656.8334
$ws.Range("K122").Value = 4970.5002
$ws.Range("M122").Value = -2520.5002
# Row 128
$ws.Range("H128").Value = 49000
$ws.Range("J128").Value = 49000
$ws.Range("L128").Value = 49000
$ws.Range("N128").Value = -58960

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 133926.94
$ws.Range("I14").Value = 133926.94
$ws.Range("K14").Value = 401780.82
$ws.Range("M14").Value = -401607.82
# Row 55
$ws.Range("H55").Value = 1908.8
$ws.Range("J55").Value = 1750
$ws.Range("L55").Value = 5250
$ws.Range("N55").Value = -5604
# Row 126
$ws.Range("H126").Value = 7007.7334
$ws.Range("I126").Value = 4261.6
$ws.Range("K126").Value = 12784.8
$ws.Range("M126").Value = -7844.800000000001
# Row 129
$ws.Range("H129").Value = 1465.8182
$ws.Range("I129").Value = 609.1
$ws.Range("K129").Value = 1827.3
$ws.Range("M129").Value = 3172.7
# Row 130
$ws.Range("H130").Value = 2686.6667
$ws.Range("I130").Value = 1825.6
$ws.Range("K130").Value = 5476.799999999999
$ws.Range("M130").Value = -456.7999999999993
# Row 131
$ws.Range("H131").Value = 20076.25
$ws.Range("J131").Value = 29689.875
$ws.Range("L131").Value = 89069.625
$ws.Range("N131").Value = -99149.625
# Row 136
$ws.Range("H136").Value = 6124.75
$ws.Range("I136").Value = 6124.75
$ws.Range("K136").Value = 18374.25
$ws.Range("M136").Value = -13274.25
# Row 137
$ws.Range("H137").Value = 5179.5293
$ws.Range("J137").Value = 6526.909
$ws.Range("L137").Value = 19580.727
$ws.Range("N137").Value = -29780.727
# Row 139
$ws.Range("H139").Value = 3984.8696
$ws.Range("I139").Value = 3000
$ws.Range("K139").Value = 9000
$ws.Range("M139").Value = -3860

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2974.875
$ws.Range("I7").Value = 2685.5715
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 2685.5715
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -2573.5715
$ws.Range("N7").Value = -5224
# Row 16
$ws.Range("H16").Value = 543.44446
$ws.Range("I16").Value = 541.8333
$ws.Range("J16").Value = 546.6667
$ws.Range("K16").Value = 541.8333
$ws.Range("L16").Value = 546.6667
$ws.Range("M16").Value = -371.8333
$ws.Range("N16").Value = -886.6667
# Row 126
$ws.Range("H126").Value = 2974.875
$ws.Range("I126").Value = 2685.5715
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 8056.7145
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -5586.7145
$ws.Range("N126").Value = -19940
# Row 132
$ws.Range("H132").Value = 806792.25
$ws.Range("I132").Value = 912383.7
$ws.Range("K132").Value = 2737151.1
$ws.Range("M132").Value = -2734621.1
# Row 140
$ws.Range("H140").Value = 98652
$ws.Range("J140").Value = 98652
$ws.Range("L140").Value = 98652
$ws.Range("N140").Value = -109012

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
# Row 102
$ws.Range("H102").Value = 100337
$ws.Range("J102").Value = 100337
$ws.Range("L102").Value = 100337
$ws.Range("N102").Value = -106827
